$wb = $excel.ActiveWorkbook

# --- "Financials" sheet: new "heading" column (K) labelling the block as
#     "Total Revenue: 2023-2024", mirroring the existing id/label/percentage/
#     title columns (C/E/G/I) ---
$wsFin = $wb.Worksheets.Item("Financials")
$wsFin.Range("K6").Value = "heading"
$wsFin.Range("K8:K16").Value = "Total Revenue: 2023-2024"
[void]$wsFin.Range("K6:K16").Select()

# --- "Demand_view" sheet: same new "heading" column (M), connected through
#     the database like the Financials sheet ---
$wsDem = $wb.Worksheets.Item("Demand_view")
$wsDem.Range("M7").Value = "heading"
$wsDem.Range("M9:M12").Value = "Total Revenue: 2023-2024"
[void]$wsDem.Range("N14").Select()
